$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 15, shifting existing rows 15-28 down to 16-29,
# to make room for the new pedestrian wait normalisation factor parameter.
$ws.Rows.Item(15).Insert()

# The freshly inserted row doesn't perfectly inherit D-column borders, so
# copy the number formatting down from the row above (same style as the
# other "value" cells in this block) before writing the new values.
$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(15, 2).Value = "ped_wait_norm"
$ws.Cells.Item(15, 3).Value = "pedestrian wait normalisation factor"
$ws.Cells.Item(15, 4).Value = 100

$ws.Range("D15").Select() | Out-Null
